# Insert a new data row at row 42 (pushing the existing rows 42-79 down to
# 43-80) and populate it with a new "Naranja" price-report record for
# Agrícola del Norte S.A. de Arica.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 42..79 down to 43..80, leaving a blank row 42 to fill in.
$ws.Rows.Item(42).Insert()

$ws.Cells.Item(42, 1).Value  = 1
$ws.Cells.Item(42, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(42, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(42, 4).Value  = 44566
$ws.Cells.Item(42, 5).Value  = 15
$ws.Cells.Item(42, 6).Value  = "Fruta"
$ws.Cells.Item(42, 7).Value  = 100102
$ws.Cells.Item(42, 8).Value  = "Cítricos"
$ws.Cells.Item(42, 9).Value  = 100102005
$ws.Cells.Item(42, 10).Value = "Naranja"
$ws.Cells.Item(42, 11).Value = "Valencia"
$ws.Cells.Item(42, 12).Value = "Tercera"
$ws.Cells.Item(42, 13).Value = 250
$ws.Cells.Item(42, 14).Value = 750
$ws.Cells.Item(42, 15).Value = 800
$ws.Cells.Item(42, 16).Value = 775
$ws.Cells.Item(42, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(42, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(42, 19).Value = 775
$ws.Cells.Item(42, 20).Value = 1
